# Generate Report for Handoff
# Updates the "Ready for handoff" rows (7,8,9,12,13,14) across the
# Overview / zh-cn / de-de sheets: stamps the new handoff/handback
# timestamps and marks the Priority column as "ht" for zh-cn & de-de.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-24 10:21:31"
}

# --- zh-cn sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-24 10:21:26"
}

# --- de-de sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-24 10:21:31"
}
